$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 42: header-like row made of single-character shared strings spelling "infoken|"
# followed by numeric cells
$ws.Range("C42").Value = "i"
$ws.Range("D42").Value = "n"
$ws.Range("E42").Value = "f"
$ws.Range("F42").Value = "o"
$ws.Range("G42").Value = "t"
$ws.Range("H42").Value = "o"
$ws.Range("I42").Value = "k"
$ws.Range("J42").Value = "e"
$ws.Range("K42").Value = "n"
$ws.Range("L42").Value = "|"
$ws.Range("M42").Value = 3
$ws.Range("N42").Value = 2
$ws.Range("O42").Value = 1
$ws.Range("P42").Value = 5
$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = 0
$ws.Range("S42").Value = 7
$ws.Range("T42").Value = 6
$ws.Range("U42").Value = 5
$ws.Range("V42").Value = 4
$ws.Range("W42").Value = 3
$ws.Range("X42").Value = 2

# Row 43: sequential numbers 0..21 from column C to X
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 2
$ws.Range("F43").Value = 3
$ws.Range("G43").Value = 4
$ws.Range("H43").Value = 5
$ws.Range("I43").Value = 6
$ws.Range("J43").Value = 7
$ws.Range("K43").Value = 8
$ws.Range("L43").Value = 9
$ws.Range("M43").Value = 10
$ws.Range("N43").Value = 11
$ws.Range("O43").Value = 12
$ws.Range("P43").Value = 13
$ws.Range("Q43").Value = 14
$ws.Range("R43").Value = 15
$ws.Range("S43").Value = 16
$ws.Range("T43").Value = 17
$ws.Range("U43").Value = 18
$ws.Range("V43").Value = 19
$ws.Range("W43").Value = 20
$ws.Range("X43").Value = 21

# Update the view to reflect the final selection/scroll position
$ws.Range("A28").Select()
$ws.Range("AD40").Select()
